$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value, exactly matching the authoritative diff.
# Cells whose new text would otherwise be auto-recognized by Excel as a number
# are written with a temporary Text number format so they stay plain text,
# matching the original inlineStr (string) cell type; the format is cleared
# again afterwards so no extra formatting is left behind on the cell.
$updates = @(
    @{ Cell = "D2"; Value = "28.084.74"; ForceText = $False }
    @{ Cell = "E2"; Value = "  -0.41%  "; ForceText = $False }
    @{ Cell = "D3"; Value = "1.874.36"; ForceText = $False }
    @{ Cell = "E3"; Value = "  -1.97%  "; ForceText = $False }
    @{ Cell = "D4"; Value = "1.003"; ForceText = $True }
    @{ Cell = "E4"; Value = "  +0.22%  "; ForceText = $False }
    @{ Cell = "D5"; Value = "313.24"; ForceText = $True }
    @{ Cell = "E5"; Value = "  -0.37%  "; ForceText = $False }
    @{ Cell = "E6"; Value = "  +0.14%  "; ForceText = $False }
    @{ Cell = "D7"; Value = "0.5047"; ForceText = $True }
    @{ Cell = "E7"; Value = "  -0.23%  "; ForceText = $False }
    @{ Cell = "D8"; Value = "0.3839"; ForceText = $True }
    @{ Cell = "D9"; Value = "0.08549"; ForceText = $True }
    @{ Cell = "E9"; Value = "  -8.49%  "; ForceText = $False }
    @{ Cell = "D10"; Value = "1.114"; ForceText = $True }
    @{ Cell = "E10"; Value = "  -2.55%  "; ForceText = $False }
    @{ Cell = "D11"; Value = "41.29"; ForceText = $True }
    @{ Cell = "E11"; Value = "  -1.73%  "; ForceText = $False }
    @{ Cell = "D12"; Value = "6.292"; ForceText = $True }
    @{ Cell = "E12"; Value = "  -1.85%  "; ForceText = $False }
    @{ Cell = "B13"; Value = "WrappedEther"; ForceText = $False }
    @{ Cell = "C13"; Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; ForceText = $False }
    @{ Cell = "D13"; Value = "1.880.21"; ForceText = $False }
    @{ Cell = "E13"; Value = "  -1.89%  "; ForceText = $False }
    @{ Cell = "B14"; Value = "Solana"; ForceText = $False }
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; ForceText = $False }
    @{ Cell = "D14"; Value = "20.64"; ForceText = $True }
    @{ Cell = "E14"; Value = "  -1.75%  "; ForceText = $False }
    @{ Cell = "D15"; Value = "7.206"; ForceText = $True }
    @{ Cell = "E15"; Value = "  -1.72%  "; ForceText = $False }
    @{ Cell = "E16"; Value = "  +0.20%  "; ForceText = $False }
    @{ Cell = "E17"; Value = "  -2.74%  "; ForceText = $False }
    @{ Cell = "D18"; Value = "90.99"; ForceText = $True }
    @{ Cell = "E18"; Value = "  -1.75%  "; ForceText = $False }
    @{ Cell = "D19"; Value = "0.06628"; ForceText = $True }
    @{ Cell = "E19"; Value = "  +0.02%  "; ForceText = $False }
    @{ Cell = "D20"; Value = "18.08"; ForceText = $True }
    @{ Cell = "E20"; Value = "  +0.31%  "; ForceText = $False }
    @{ Cell = "E21"; Value = "  +0.20%  "; ForceText = $False }
    @{ Cell = "D22"; Value = "6.093"; ForceText = $True }
    @{ Cell = "E22"; Value = "  -2.30%  "; ForceText = $False }
    @{ Cell = "D23"; Value = "28.124.70"; ForceText = $False }
    @{ Cell = "E23"; Value = "  -0.47%  "; ForceText = $False }
    @{ Cell = "D24"; Value = "11.37"; ForceText = $True }
    @{ Cell = "E24"; Value = "  -1.64%  "; ForceText = $False }
    @{ Cell = "D25"; Value = "2.267"; ForceText = $True }
    @{ Cell = "E25"; Value = "  -2.71%  "; ForceText = $False }
    @{ Cell = "D26"; Value = "2.585"; ForceText = $True }
    @{ Cell = "E26"; Value = "  -0.41%  "; ForceText = $False }
    @{ Cell = "D27"; Value = "2.091.53"; ForceText = $False }
    @{ Cell = "E27"; Value = "  -2.14%  "; ForceText = $False }
    @{ Cell = "D28"; Value = "20.71"; ForceText = $True }
    @{ Cell = "E28"; Value = "  -2.36%  "; ForceText = $False }
    @{ Cell = "D29"; Value = "156.88"; ForceText = $True }
    @{ Cell = "E29"; Value = "  -0.93%  "; ForceText = $False }
    @{ Cell = "D30"; Value = "126.26"; ForceText = $True }
    @{ Cell = "E30"; Value = "  -0.85%  "; ForceText = $False }
    @{ Cell = "E31"; Value = "  -1.43%  "; ForceText = $False }
    @{ Cell = "E32"; Value = "  -4.25%  "; ForceText = $False }
    @{ Cell = "E33"; Value = "  -1.01%  "; ForceText = $False }
    @{ Cell = "D34"; Value = "3.585"; ForceText = $True }
    @{ Cell = "E34"; Value = "  -0.78%  "; ForceText = $False }
    @{ Cell = "D35"; Value = "9.622"; ForceText = $True }
    @{ Cell = "E35"; Value = "  -0.79%  "; ForceText = $False }
    @{ Cell = "D36"; Value = "0.02450"; ForceText = $True }
    @{ Cell = "E36"; Value = "  +0.36%  "; ForceText = $False }
    @{ Cell = "D37"; Value = "0.06582"; ForceText = $True }
    @{ Cell = "E37"; Value = "  -1.93%  "; ForceText = $False }
    @{ Cell = "D38"; Value = "0.2180"; ForceText = $True }
    @{ Cell = "E38"; Value = "  -1.70%  "; ForceText = $False }
    @{ Cell = "D39"; Value = "1.213"; ForceText = $True }
    @{ Cell = "E39"; Value = "  -2.80%  "; ForceText = $False }
    @{ Cell = "D40"; Value = "1.242"; ForceText = $True }
    @{ Cell = "E40"; Value = "  -3.07%  "; ForceText = $False }
    @{ Cell = "D41"; Value = "0.6376"; ForceText = $True }
    @{ Cell = "E41"; Value = "  -2.69%  "; ForceText = $False }
    @{ Cell = "D42"; Value = "11.43"; ForceText = $True }
    @{ Cell = "E42"; Value = "  -1.17%  "; ForceText = $False }
    @{ Cell = "D43"; Value = "4.896"; ForceText = $True }
    @{ Cell = "E43"; Value = "  -2.72%  "; ForceText = $False }
    @{ Cell = "D44"; Value = "13.18"; ForceText = $True }
    @{ Cell = "E44"; Value = "  -1.42%  "; ForceText = $False }
    @{ Cell = "D45"; Value = "0.6004"; ForceText = $True }
    @{ Cell = "E45"; Value = "  -2.25%  "; ForceText = $False }
    @{ Cell = "D46"; Value = "1.284"; ForceText = $True }
    @{ Cell = "E46"; Value = "  -1.40%  "; ForceText = $False }
    @{ Cell = "D47"; Value = "3.673"; ForceText = $True }
    @{ Cell = "E47"; Value = "  -1.39%  "; ForceText = $False }
    @{ Cell = "B48"; Value = "NEARProtocol"; ForceText = $False }
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; ForceText = $False }
    @{ Cell = "D48"; Value = "1.989"; ForceText = $True }
    @{ Cell = "E48"; Value = "  -2.19%  "; ForceText = $False }
    @{ Cell = "B49"; Value = "EOS"; ForceText = $False }
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"; ForceText = $False }
    @{ Cell = "D49"; Value = "1.223"; ForceText = $True }
    @{ Cell = "E49"; Value = "  +2.65%  "; ForceText = $False }
    @{ Cell = "D50"; Value = "121.41"; ForceText = $True }
    @{ Cell = "E50"; Value = "  -0.94%  "; ForceText = $False }
    @{ Cell = "D51"; Value = "80.63"; ForceText = $True }
    @{ Cell = "E51"; Value = "  +2.74%  "; ForceText = $False }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.ClearFormats()
    } else {
        $rng.Value = $u.Value
    }
}
